{"js": "// Replace the date header and each two-digit \u00f7 one-digit problem in the\n// practice-sheet table. Every source string below occurs exactly once in\n// the document, so a body-wide literal `search()` + `insertText('Replace')`\n// per pair is safe even though some replacement texts equal OTHER pairs'\n// source texts (e.g. \"96\u00f74=\" is both a source and, elsewhere, a target) \u2014\n// each search/replace is resolved against the ORIGINAL text before any of\n// these edits run, so there is no cross-talk between pairs.\nconst replacements = [\n  [\"2024-05-31 Friday\", \"2024-06-01 Saturday\"],\n  [\"99\u00f79=\", \"27\u00f73=\"],\n  [\"96\u00f74=\", \"58\u00f72=\"],\n  [\"52\u00f78=\", \"22\u00f73=\"],\n  [\"91\u00f73=\", \"50\u00f76=\"],\n  [\"74\u00f73=\", \"16\u00f74=\"],\n  [\"84\u00f77=\", \"74\u00f79=\"],\n  [\"69\u00f77=\", \"24\u00f76=\"],\n  [\"41\u00f72=\", \"66\u00f77=\"],\n  [\"20\u00f75=\", \"19\u00f79=\"],\n  [\"79\u00f72=\", \"39\u00f73=\"],\n  [\"36\u00f72=\", \"24\u00f73=\"],\n  [\"92\u00f72=\", \"55\u00f75=\"],\n  [\"38\u00f73=\", \"77\u00f79=\"],\n  [\"37\u00f75=\", \"48\u00f75=\"],\n  [\"58\u00f72=\", \"99\u00f72=\"],\n  [\"68\u00f74=\", \"24\u00f78=\"],\n  [\"33\u00f76=\", \"26\u00f77=\"],\n  [\"76\u00f73=\", \"96\u00f74=\"],\n  [\"63\u00f75=\", \"62\u00f74=\"],\n  [\"97\u00f79=\", \"95\u00f78=\"],\n  [\"54\u00f73=\", \"29\u00f74=\"],\n  [\"90\u00f74=\", \"68\u00f78=\"],\n  [\"40\u00f79=\", \"85\u00f75=\"],\n  [\"14\u00f73=\", \"11\u00f74=\"],\n  [\"52\u00f73=\", \"79\u00f74=\"],\n];\n\nconst body = context.document.body;\n\n// Resolve every search BEFORE mutating anything so later inserts can't\n// shadow an earlier pair's source text (see note above).\nconst hits = replacements.map(([from]) => body.search(from, { matchCase: true }));\nhits.forEach((h) => h.load(\"items\"));\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [, to] = replacements[i];\n  const items = hits[i].items;\n  for (const range of items) {\n    range.insertText(to, \"Replace\");\n  }\n}\nawait context.sync();\n", "ps1": "# Update the worksheet date header and the 25 two-digit / one-digit\n# division problems laid out 5-per-row in the single table, one problem\n# row every 4th table row (rows 1, 5, 9, 13, 17; the rows in between are\n# blank answer rows).\n#\n# Each Find/Replace is scoped to ONE table cell by selecting that cell's\n# Range and running Find against the Selection. This matters because\n# several replacement strings equal ANOTHER cell's original text (e.g.\n# row1/col2 \"96\u00f74=\" becomes \"58\u00f72=\" while row9/col5 starts as \"58\u00f72=\" and\n# becomes \"99\u00f72=\"; row13/col3 \"76\u00f73=\" becomes \"96\u00f74=\"). Scoping the\n# search+replace to each cell's own selection keeps every replacement\n# targeted at exactly the one cell it belongs to, with no cross-talk\n# between cells that momentarily share the same text mid-script.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\nfunction Replace-InCell($row, $col, $findText, $replaceText) {\n    $cell = $tbl.Cell($row, $col)\n    $cell.Range.Select()\n    $sel = $word.Selection\n    $sel.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 0, $false, $replaceText, 1) | Out-Null\n}\n\n$problemRows = @(1, 5, 9, 13, 17)\n$grid = @(\n    @(\"99\u00f79=\", \"27\u00f73=\"), @(\"96\u00f74=\", \"58\u00f72=\"), @(\"52\u00f78=\", \"22\u00f73=\"), @(\"91\u00f73=\", \"50\u00f76=\"), @(\"74\u00f73=\", \"16\u00f74=\"),\n    @(\"84\u00f77=\", \"74\u00f79=\"), @(\"69\u00f77=\", \"24\u00f76=\"), @(\"41\u00f72=\", \"66\u00f77=\"), @(\"20\u00f75=\", \"19\u00f79=\"), @(\"79\u00f72=\", \"39\u00f73=\"),\n    @(\"36\u00f72=\", \"24\u00f73=\"), @(\"92\u00f72=\", \"55\u00f75=\"), @(\"38\u00f73=\", \"77\u00f79=\"), @(\"37\u00f75=\", \"48\u00f75=\"), @(\"58\u00f72=\", \"99\u00f72=\"),\n    @(\"68\u00f74=\", \"24\u00f78=\"), @(\"33\u00f76=\", \"26\u00f77=\"), @(\"76\u00f73=\", \"96\u00f74=\"), @(\"63\u00f75=\", \"62\u00f74=\"), @(\"97\u00f79=\", \"95\u00f78=\"),\n    @(\"54\u00f73=\", \"29\u00f74=\"), @(\"90\u00f74=\", \"68\u00f78=\"), @(\"40\u00f79=\", \"85\u00f75=\"), @(\"14\u00f73=\", \"11\u00f74=\"), @(\"52\u00f73=\", \"79\u00f74=\")\n)\n\nfor ($i = 0; $i -lt $grid.Count; $i++) {\n    $row = $problemRows[[Math]::Floor($i / 5)]\n    $col = ($i % 5) + 1\n    Replace-InCell $row $col $grid[$i][0] $grid[$i][1]\n}\n\n# Date header \u2014 first paragraph in the body, outside the table. This text\n# is unique in the document, so a plain whole-document Find/Replace (not\n# Selection-scoped \u2014 selecting table cells above can leave\n# Document.Paragraphs(1) pointing at a stale/reseated paragraph) is safe\n# here and avoids that quirk entirely.\n$d.Content.Find.Execute(\"2024-05-31 Friday\", $false, $false, $false, $false, $false, $true, 0, $false, \"2024-06-01 Saturday\", 1) | Out-Null\n"}
